$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the existing last-row key cell (A82) the "vertical-top" alignment style
# (numFmtId 164 / General, vertical=top) -- this creates a new cellXfs entry.
$ws.Range("A82").VerticalAlignment = -4160

# Append the new product row.
$ws.Range("A83").Value = 12583189
$ws.Range("B83").Value = 105

# Move the view so the newly-added row is visible and selected, matching
# the author's final cursor position.
$ws.Range("B83").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 74
$win.ScrollColumn = 1
